# Apply the "May 9th" data refresh for the walkingToRunning sample sheet.
# - Columns C:H (ax..gz) for existing rows 2-21 are replaced with a later
#   slice of the same sensor stream (effectively the window advanced by 12
#   samples).
# - 10 new rows (22-31) are appended, continuing the timestamp sequence
#   (2000..2900) with the "walkingToRunning" label.
# - The sheet dimension grows from A1:H21 to A1:H31 automatically once the
#   new rows are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 30,8

# row 2: timestamp=0
$data[0,0] = 0
$data[0,1] = "walkingToRunning"
$data[0,2] = 3.43527889251709
$data[0,3] = -5.793748378753662
$data[0,4] = 33.75442886352539
$data[0,5] = 1.886958264387571
$data[0,6] = -2.518787022737357
$data[0,7] = -1.387751776667747
# row 3: timestamp=100
$data[1,0] = 100
$data[1,1] = "walkingToRunning"
$data[1,2] = 3.264275550842285
$data[1,3] = -8.339526176452637
$data[1,4] = 33.21075820922852
$data[1,5] = -1.408535764767576
$data[1,6] = -0.4590572393857499
$data[1,7] = 1.789505992944425
# row 4: timestamp=200
$data[2,0] = 200
$data[2,1] = "walkingToRunning"
$data[2,2] = -2.304214000701904
$data[2,3] = -44.59181213378906
$data[2,4] = 10.13076019287109
$data[2,5] = -9.029990652891366
$data[2,6] = -0.2165457317462454
$data[2,7] = 2.049246124120849
# row 5: timestamp=300
$data[3,0] = 300
$data[3,1] = "walkingToRunning"
$data[3,2] = -2.632879257202148
$data[3,3] = -47.35254669189453
$data[3,4] = 16.20354652404785
$data[3,5] = -1.633290004730144
$data[3,6] = 0.8114918470382817
$data[3,7] = 0.4343748807907029
# row 6: timestamp=400
$data[4,0] = 400
$data[4,1] = "walkingToRunning"
$data[4,2] = -20.92312431335449
$data[4,3] = -10.33131408691406
$data[4,4] = -20.38119125366211
$data[4,5] = 4.732972145080568
$data[4,6] = -2.583366572856895
$data[4,7] = -0.3213605839472549
# row 7: timestamp=500
$data[5,0] = 500
$data[5,1] = "walkingToRunning"
$data[5,2] = -20.04123497009277
$data[5,3] = -3.740296125411987
$data[5,4] = -22.7781810760498
$data[5,5] = 4.758649143805877
$data[5,6] = -1.287177375646699
$data[5,7] = 0.682966625690469
# row 8: timestamp=600
$data[6,0] = 600
$data[6,1] = "walkingToRunning"
$data[6,2] = -13.83809471130371
$data[6,3] = -1.457606315612793
$data[6,4] = 10.62628936767578
$data[6,5] = 1.202222677377485
$data[6,6] = 1.668132488544167
$data[6,7] = 1.826421522177189
# row 9: timestamp=700
$data[7,0] = 700
$data[7,1] = "walkingToRunning"
$data[7,2] = -15.24413108825684
$data[7,3] = -5.868594169616699
$data[7,4] = 10.52801132202148
$data[7,5] = -3.896952409010684
$data[7,6] = 0.7400126088123988
$data[7,7] = 1.755212721457845
# row 10: timestamp=800
$data[8,0] = 800
$data[8,1] = "walkingToRunning"
$data[8,2] = -68.24240875244141
$data[8,3] = -37.54502487182617
$data[8,4] = -25.38205528259277
$data[8,5] = -3.904241349146918
$data[8,6] = 8.71459445529258
$data[8,7] = 1.42584520486685
# row 11: timestamp=900
$data[9,0] = 900
$data[9,1] = "walkingToRunning"
$data[9,2] = -69.99400329589844
$data[9,3] = -36.24868774414063
$data[9,4] = -23.62584495544434
$data[9,5] = -1.888838914724495
$data[9,6] = 4.313753348130421
$data[9,7] = 0.935868626374464
# row 12: timestamp=1000
$data[10,0] = 1000
$data[10,1] = "walkingToRunning"
$data[10,2] = 3.574422836303711
$data[10,3] = -10.55799293518066
$data[10,4] = -1.943589687347412
$data[10,5] = 4.815394823367857
$data[10,6] = -4.399581175584057
$data[10,7] = -2.032744779036629
# row 13: timestamp=1100
$data[11,0] = 1100
$data[11,1] = "walkingToRunning"
$data[11,2] = 1.990053176879883
$data[11,3] = -4.801568031311035
$data[11,4] = -0.5777735710144043
$data[11,5] = 5.744059562683067
$data[11,6] = -2.208836053426431
$data[11,7] = -2.385511210331534
# row 14: timestamp=1200
$data[12,0] = 1200
$data[12,1] = "walkingToRunning"
$data[12,2] = 15.68469429016113
$data[12,3] = 14.39326095581055
$data[12,4] = 3.627434253692627
$data[12,5] = -0.9773453107246968
$data[12,6] = -5.042484849920636
$data[12,7] = -1.550856040532772
# row 15: timestamp=1300
$data[13,0] = 1300
$data[13,1] = "walkingToRunning"
$data[13,2] = 15.73915863037109
$data[13,3] = 7.167607307434082
$data[13,4] = 5.527215480804443
$data[13,5] = -5.169222395236689
$data[13,6] = -6.381682007129349
$data[13,7] = -6.139806142220171
# row 16: timestamp=1400
$data[14,0] = 1400
$data[14,1] = "walkingToRunning"
$data[14,2] = -0.1553750038146972
$data[14,3] = 10.89957809448242
$data[14,4] = -15.06903266906738
$data[14,5] = -4.023792711588003
$data[14,6] = 7.085793306277576
$data[14,7] = -4.32150216927884
# row 17: timestamp=1500
$data[15,0] = 1500
$data[15,1] = "walkingToRunning"
$data[15,2] = -0.4681458473205566
$data[15,3] = 11.8446626663208
$data[15,4] = -15.58957099914551
$data[15,5] = 1.662728214263952
$data[15,6] = 4.684056318723108
$data[15,7] = -1.641669896932774
# row 18: timestamp=1600
$data[16,0] = 1600
$data[16,1] = "walkingToRunning"
$data[16,2] = 9.029394149780272
$data[16,3] = -76.68982696533203
$data[16,4] = 2.083036422729492
$data[16,5] = 2.464888723567128
$data[16,6] = -12.57504773139953
$data[16,7] = -3.376289129257202
# row 19: timestamp=1700
$data[17,0] = 1700
$data[17,1] = "walkingToRunning"
$data[17,2] = 9.611247062683104
$data[17,3] = -76.88915252685547
$data[17,4] = 1.387777328491211
$data[17,5] = -0.0004609318306818891
$data[17,6] = -8.720019648625392
$data[17,7] = -0.9112097813533104
# row 20: timestamp=1800
$data[18,0] = 1800
$data[18,1] = "walkingToRunning"
$data[18,2] = -11.94809913635254
$data[18,3] = 10.77702140808106
$data[18,4] = 4.135905742645264
$data[18,5] = -6.561141893869462
$data[18,6] = 1.067614199106516
$data[18,7] = 5.132379430073989
# row 21: timestamp=1900
$data[19,0] = 1900
$data[19,1] = "walkingToRunning"
$data[19,2] = -66.71106719970703
$data[19,3] = -81.49526214599609
$data[19,4] = 59.02106475830078
$data[19,5] = -3.78096956106336
$data[19,6] = -6.149249487656785
$data[19,7] = 3.56855486356297
# row 22: timestamp=2000
$data[20,0] = 2000
$data[20,1] = "walkingToRunning"
$data[20,2] = -43.36343383789063
$data[20,3] = -8.654863357543945
$data[20,4] = -11.20723152160644
$data[20,5] = 4.143512133451609
$data[20,6] = -9.295757961273193
$data[20,7] = 1.042706482227031
# row 23: timestamp=2100
$data[21,0] = 2100
$data[21,1] = "walkingToRunning"
$data[21,2] = -40.74956512451172
$data[21,3] = -3.842054843902588
$data[21,4] = -15.99789905548096
$data[21,5] = 3.805606016746037
$data[21,6] = 1.494310085590225
$data[21,7] = 2.446077621900006
# row 24: timestamp=2200
$data[22,0] = 2200
$data[22,1] = "walkingToRunning"
$data[22,2] = -8.164802551269531
$data[22,3] = 7.02423906326294
$data[22,4] = -4.845988273620605
$data[22,5] = 3.119575682053099
$data[22,6] = 4.113574073864802
$data[22,7] = -0.2212830832371099
# row 25: timestamp=2300
$data[23,0] = 2300
$data[23,1] = "walkingToRunning"
$data[23,2] = -8.398514747619629
$data[23,3] = 4.322819232940674
$data[23,4] = -5.664791107177734
$data[23,5] = -2.658262938719541
$data[23,6] = 2.234168899976305
$data[23,7] = 2.193659129509576
# row 26: timestamp=2400
$data[24,0] = 2400
$data[24,1] = "walkingToRunning"
$data[24,2] = -79.12488555908203
$data[24,3] = -28.08302688598633
$data[24,4] = -31.37481689453125
$data[24,5] = -2.523965861247135
$data[24,6] = 9.674394070185159
$data[24,7] = 1.710836189526775
# row 27: timestamp=2500
$data[25,0] = 2500
$data[25,1] = "walkingToRunning"
$data[25,2] = -81.6524658203125
$data[25,3] = -26.21007919311523
$data[25,4] = -27.36414909362793
$data[25,5] = -0.4347282877335233
$data[25,6] = 7.832854989858879
$data[25,7] = -0.6348883940623381
# row 28: timestamp=2600
$data[26,0] = 2600
$data[26,1] = "walkingToRunning"
$data[26,2] = -9.369277954101562
$data[26,3] = -24.22552871704102
$data[26,4] = -6.347414493560791
$data[26,5] = 2.215493936378243
$data[26,6] = -6.333109965691214
$data[26,7] = -3.700846394667294
# row 29: timestamp=2700
$data[27,0] = 2700
$data[27,1] = "walkingToRunning"
$data[27,2] = -8.832977294921875
$data[27,3] = -17.93884468078613
$data[27,4] = -7.180578708648682
$data[27,5] = 3.503521535946755
$data[27,6] = -4.84552946755515
$data[27,7] = -5.782033076653083
# row 30: timestamp=2800
$data[28,0] = 2800
$data[28,1] = "walkingToRunning"
$data[28,2] = -12.03362846374512
$data[28,3] = 20.78299713134766
$data[28,4] = 6.19299840927124
$data[28,5] = 0.5721309762734461
$data[28,6] = -0.04020348672683609
$data[28,7] = -0.7832388446880931
# row 31: timestamp=2900
$data[29,0] = 2900
$data[29,1] = "walkingToRunning"
$data[29,2] = -11.60475063323975
$data[29,3] = 14.68077850341797
$data[29,4] = 6.681206226348877
$data[29,5] = -4.998460626602359
$data[29,6] = -7.086485686898527
$data[29,7] = 0.6293182730674433

# Single bulk write mirrors how Excel applies a pasted block; the sheet's
# <dimension> is recomputed by the engine on save.
$ws.Range("A2:H31").Value = $data
